$wb = $excel.ActiveWorkbook

# Update "想去人数" (attendee count) figures on both the "展览" and "全部类型" sheets.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8567
    $ws.Range("F4").Value = 388
    $ws.Range("F5").Value = 29
}
